# Apply the "Add a working example data-sheet" edit to Blad1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")
$ws.Activate()

# --- Row 2 (Bovenwaarde) ---
$ws.Range("I2").Value = 15
$ws.Range("J2").Value = 165
$ws.Range("L2").Value = 380
$ws.Range("P2").Value = 100
$ws.Range("R2").Value = 1
$ws.Range("S2").Value = "inf"
$ws.Range("V2").Formula = "=2/19"
$ws.Range("W2").Formula = "=2/19"

# --- Row 3 (Onderwaarde) ---
$ws.Range("F3").Value = 950
$ws.Range("G3").Value = 550
$ws.Range("H3").Value = 85
$ws.Range("J3").Value = 148
$ws.Range("L3").Value = 300
$ws.Range("N3").Value = 100
$ws.Range("P3").Value = 20
$ws.Range("Q3").Value = 20
$ws.Range("R3").Value = 1
$ws.Range("S3").Formula = "=13/19"
$ws.Range("T3").Formula = "=4/19"
$ws.Range("U3").Formula = "=4/19"

# --- Row 5 (Snijmais) ---
$ws.Range("T5").Value = 1

# --- Row 6 (Graskuil) ---
$ws.Range("U6").Value = 1

# --- Row 18 (Sojaschroot nonGMO 44%) ---
$ws.Range("W18").Value = 1

# --- Row 24: new formatted (date, "d-mmm") but otherwise empty cell ---
$ws.Range("B24").NumberFormat = "d-mmm"

# --- Window scroll position / selection ---
$ws.Range("W19").Select()
